{"js": "// Minor changes in Prompts:\n// 1) \" find my various arch\" -> \" find various arch\"\n// 2) \" efficiently and accurately and can be run on colab irrespective of any issues.\" -> \" efficiently and accurately.\"\n// 3) \"halluconated\" -> \"hallucinated\" (typo fix)\n\nconst body = context.document.body;\n\n// 1) Remove \"my \" before \"various arch\" (first paragraph).\nlet results = body.search(\" find my various arch\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\" find various arch\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Drop \"and can be run on colab irrespective of any issues\" \u2014 keep the\n//    sentence ending at \"efficiently and accurately.\"\nresults = body.search(\n  \" efficiently and accurately and can be run on colab irrespective of any issues.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\" efficiently and accurately.\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Fix the misspelling \"halluconated\" -> \"hallucinated\".\nresults = body.search(\"halluconated\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"hallucinated\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Minor changes in Prompts:\n# 1) \" find my various arch\" -> \" find various arch\"\n# 2) \" efficiently and accurately and can be run on colab irrespective of any issues.\" -> \" efficiently and accurately.\"\n# 3) \"halluconated\" -> \"hallucinated\" (typo fix)\n\n$d = $word.ActiveDocument\n\n# 1) Remove \"my \" before \"various arch\" (first paragraph).\n$find1 = $d.Content.Find\n$find1.Execute(\" find my various arch\", $false, $false, $false, $false, $false, $true, 1, $false, \" find various arch\", 2)\n\n# 2) Drop \"and can be run on colab irrespective of any issues\" \u2014 keep the\n#    sentence ending at \"efficiently and accurately.\"\n$find2 = $d.Content.Find\n$find2.Execute(\" efficiently and accurately and can be run on colab irrespective of any issues.\", $false, $false, $false, $false, $false, $true, 1, $false, \" efficiently and accurately.\", 2)\n\n# 3) Fix the misspelling \"halluconated\" -> \"hallucinated\".\n$find3 = $d.Content.Find\n$find3.Execute(\"halluconated\", $false, $false, $false, $false, $false, $true, 1, $false, \"hallucinated\", 2)\n"}
